$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.712.51'
$ws.Range("E2").Value = '  +4.50%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.762.80'
$ws.Range("E3").Value = '  +5.06%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '116.57'
$ws.Range("E5").Value = '  +3.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '333.22'
$ws.Range("E6").Value = '  +2.88%  '
$ws.Range("E7").Value = '  +2.16%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +6.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.87'
$ws.Range("E10").Value = '  +4.76%  '
$ws.Range("E11").Value = '  +5.81%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.19'
$ws.Range("E12").Value = '  +1.87%  '
$ws.Range("E13").Value = '  +1.93%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.67'
$ws.Range("E14").Value = '  +5.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.195.22'
$ws.Range("E15").Value = '  +5.17%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.767.66'
$ws.Range("E16").Value = '  +5.36%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.891'
$ws.Range("E17").Value = '  +3.47%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.715.82'
$ws.Range("E18").Value = '  +4.75%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.20'
$ws.Range("E19").Value = '  +5.60%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.54'
$ws.Range("E20").Value = '  +4.55%  '
$ws.Range("E21").Value = '  +2.49%  '
$ws.Range("E22").Value = '  +2.89%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '278.54'
$ws.Range("E23").Value = '  +3.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.63'
$ws.Range("E24").Value = '  +1.40%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.67'
$ws.Range("E25").Value = '  +5.45%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.83'
$ws.Range("E26").Value = '  +2.27%  '
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.19'
$ws.Range("E28").Value = '  -1.46%  '
$ws.Range("E29").Value = '  +0.10%  '
$ws.Range("E30").Value = '  +2.67%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.12'
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '50.11'
$ws.Range("E32").Value = '  +1.08%  '
$ws.Range("E33").Value = '  +1.51%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0825'
$ws.Range("E34").Value = '  +1.34%  '
$ws.Range("E35").Value = '  -0.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '19.02'
$ws.Range("E36").Value = '  -0.13%  '
$ws.Range("E37").Value = '  +1.85%  '
$ws.Range("E38").Value = '  +1.86%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.24'
$ws.Range("E39").Value = '  +3.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0352'
$ws.Range("E40").Value = '  +9.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '127.13'
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '23.22'
$ws.Range("E42").Value = '  +3.68%  '
$ws.Range("E43").Value = '  +3.25%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.30'
$ws.Range("E44").Value = '  +7.42%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.45'
$ws.Range("E45").Value = '  +13.58%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.091.58'
$ws.Range("E46").Value = '  +1.46%  '
$ws.Range("E47").Value = '  +2.52%  '
$ws.Range("E48").Value = '  +4.75%  '
$ws.Range("E49").Value = '  +6.18%  '
$ws.Range("E50").Value = '  +1.06%  '
$ws.Range("E51").Value = '  +1.49%  '
